$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "iaest-measure:sector-actividad"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "iaest-measure:sexo"

# Row 3 updates
$ws.Range("B3").Value = "medida"
$ws.Range("F3").Value = "dim"
$ws.Range("I3").Value = "medida"

# Row 4 updates
$ws.Range("B4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Municipio"
$ws.Range("H4").Value = "URI-Comunidad"
$ws.Range("I4").Value = "xsd:int"

# Row 5 removed entirely
$ws.Rows.Item(5).Delete()
